$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 40
$ws.Cells.Item($row, 1).Value = 44685
$ws.Cells.Item($row, 2).Value = "BUAM"
$ws.Cells.Item($row, 3).Value = 121
$ws.Cells.Item($row, 4).Value = "Sherbrooke"
$ws.Cells.Item($row, 5).Value = "Estrie"
$ws.Cells.Item($row, 6).Value = "A"
$ws.Cells.Item($row, 7).Value = "Min. Cote 1"
$ws.Cells.Item($row, 9).Value = "Oreste Loica"
